$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full A1:F16 table (header + 15 data rows) as a single 2D array
# to match the new layout: Attribute, Y, R², RMSE, Offset, Slope
$data = New-Object 'object[,]' 16,6

$data[0,0] = "Attribute"
$data[0,1] = "Y"
$data[0,2] = "R²"
$data[0,3] = "RMSE"
$data[0,4] = "Offset"
$data[0,5] = "Slope"

$data[1,0] = "SST"
$data[1,1] = "Referência"
$data[1,2] = 0.8051743426048514
$data[1,3] = 1.148837929164844
$data[1,4] = 2.711639164099223
$data[1,5] = 0.8051743426048512

$data[2,0] = "SST"
$data[2,1] = "Predição"
$data[2,2] = 0.7383382171879158
$data[2,3] = 1.331391679010095
$data[2,4] = 3.067603339747406
$data[2,5] = 0.7793988667851736

$data[3,0] = "SST"
$data[3,1] = "Validação"
$data[3,2] = 0.7653390937403943
$data[3,3] = 1.088270631842013
$data[3,4] = 1.31771668119452
$data[3,5] = 0.9160208818725738

$data[4,0] = "PH"
$data[4,1] = "Referência"
$data[4,2] = 0.6243669284831987
$data[4,3] = 0.1919195916110745
$data[4,4] = 1.240089264609229
$data[4,5] = 0.6243669284832003

$data[5,0] = "PH"
$data[5,1] = "Predição"
$data[5,2] = 0.3990234459236567
$data[5,3] = 0.2427538713902996
$data[5,4] = 1.592154324702855
$data[5,5] = 0.5180921853576327

$data[6,0] = "PH"
$data[6,1] = "Validação"
$data[6,2] = 0.4840792817312685
$data[6,3] = 0.177323052565511
$data[6,4] = 0.8219452083264076
$data[6,5] = 0.7495504091221834

$data[7,0] = "AT"
$data[7,1] = "Referência"
$data[7,2] = 0.494665790894792
$data[7,3] = 0.3975020114248576
$data[7,4] = 0.5687060065510533
$data[7,5] = 0.494665790894792

$data[8,0] = "AT"
$data[8,1] = "Predição"
$data[8,2] = 0.3953592680174114
$data[8,3] = 0.4348091564766918
$data[8,4] = 0.6175741326214125
$data[8,5] = 0.4532763022904974

$data[9,0] = "AT"
$data[9,1] = "Validação"
$data[9,2] = 0.6378037674957653
$data[9,3] = 0.3090896363122474
$data[9,4] = -0.2282144515090805
$data[9,5] = 1.154882760098361

$data[10,0] = "FIRMEZA (N)"
$data[10,1] = "Referência"
$data[10,2] = 0.6127065482896518
$data[10,3] = 62.91116033880429
$data[10,4] = 199.1192575734262
$data[10,5] = 0.6127065482896522

$data[11,0] = "FIRMEZA (N)"
$data[11,1] = "Predição"
$data[11,2] = 0.3789994062683451
$data[11,3] = 79.66244398043682
$data[11,4] = 258.2357080583688
$data[11,5] = 0.4980083127535585

$data[12,0] = "FIRMEZA (N)"
$data[12,1] = "Validação"
$data[12,2] = 0.5474897551728245
$data[12,3] = 55.06736389728003
$data[12,4] = 49.75257411922672
$data[12,5] = 0.8799175796427521

$data[13,0] = "UBS (%)"
$data[13,1] = "Referência"
$data[13,2] = 0.6726734889966415
$data[13,3] = 1.800229857848642
$data[13,4] = 5.000182059074024
$data[13,5] = 0.6726734889966427

$data[14,0] = "UBS (%)"
$data[14,1] = "Predição"
$data[14,2] = 0.597283287976039
$data[14,3] = 1.996812083006219
$data[14,4] = 5.443534267243627
$data[14,5] = 0.642703459905562

$data[15,0] = "UBS (%)"
$data[15,1] = "Validação"
$data[15,2] = 0.7440479612605038
$data[15,3] = 1.213024940458897
$data[15,4] = 1.730532550889199
$data[15,5] = 0.8895905442856631

$ws.Range("A1:F16").Value = $data

